# Updated cryptos list on Tue Aug  8 07:24:52 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Force the cell to keep a literal text value (the source price strings
    # look numeric, e.g. "0.9986" or "243.10", and Excel's COM Value setter
    # would otherwise coerce them into real numbers). Using a temporary
    # text number-format preserves the string, then ClearFormats drops the
    # now-unneeded formatting so the cell style matches the original
    # (un-styled) cell.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).ClearFormats()
}

# Row 2 - Bitcoin
Set-TextCell "D2" "29.185.34"
$ws.Range("E2").Value = "  +0.30%  "

# Row 3 - Ethereum
Set-TextCell "D3" "1.830.73"
$ws.Range("E3").Value = "  -0.26%  "

# Row 4 - TetherUSD
Set-TextCell "D4" "0.9986"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5 - BNB
Set-TextCell "D5" "242.81"
$ws.Range("E5").Value = "  -0.20%  "

# Row 6 - XRP
Set-TextCell "D6" "0.6194"
$ws.Range("E6").Value = "  +0.43%  "

# Row 7 - USDC
Set-TextCell "D7" "0.9999"

# Row 8 - Dogecoin
Set-TextCell "D8" "0.07348"
$ws.Range("E8").Value = "  -1.56%  "

# Row 9 - Cardano
Set-TextCell "D9" "0.2902"
$ws.Range("E9").Value = "  -0.63%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  +0.29%  "

# Row 11 - TRON
Set-TextCell "D11" "0.07653"
$ws.Range("E11").Value = "  -0.53%  "

# Row 12 - WrappedEther
Set-TextCell "D12" "1.841.00"
$ws.Range("E12").Value = "  +0.20%  "

# Row 13 - Polkadot
Set-TextCell "D13" "4.968"
$ws.Range("E13").Value = "  -0.79%  "

# Row 14 - Polygon
Set-TextCell "D14" "0.6682"
$ws.Range("E14").Value = "  -0.54%  "

# Row 15 - Litecoin
Set-TextCell "D15" "82.46"
$ws.Range("E15").Value = "  -0.19%  "

# Row 16 - ShibaInu
Set-TextCell "D16" "0.000008954"
$ws.Range("E16").Value = "  -3.36%  "

# Row 17 - Uniswap
Set-TextCell "D17" "5.844"
$ws.Range("E17").Value = "  -1.42%  "

# Row 18 - WrappedBTC
Set-TextCell "D18" "29.177.66"
$ws.Range("E18").Value = "  +0.32%  "

# Row 19 - WrappedliquidstakedEther2.0
Set-TextCell "D19" "2.084.76"
$ws.Range("E19").Value = "  -0.09%  "

# Row 20 - BitcoinCash
Set-TextCell "D20" "235.00"
$ws.Range("E20").Value = "  +1.74%  "

# Row 21 - Avalanche
$ws.Range("E21").Value = "  -1.39%  "

# Row 22 - Dai
Set-TextCell "D22" "0.9998"
$ws.Range("E22").Value = "  -0.22%  "

# Row 23 - Chainlink
Set-TextCell "D23" "7.352"
$ws.Range("E23").Value = "  +2.40%  "

# Row 24 - BinanceUSD
Set-TextCell "D24" "0.9992"
$ws.Range("E24").Value = "  -0.24%  "

# Row 25 - Monero
Set-TextCell "D25" "158.03"
$ws.Range("E25").Value = "  -1.44%  "

# Row 26 - Stellar
Set-TextCell "D26" "0.1392"
$ws.Range("E26").Value = "  +0.41%  "

# Row 27 - Cosmos
Set-TextCell "D27" "8.539"
$ws.Range("E27").Value = "  +0.36%  "

# Row 28 - EthereumClassic
Set-TextCell "D28" "17.61"
$ws.Range("E28").Value = "  -1.13%  "

# Row 29 - PancakeSwap
Set-TextCell "D29" "1.488"
$ws.Range("E29").Value = "  -0.59%  "

# Row 30 - Hedera
Set-TextCell "D30" "0.05830"
$ws.Range("E30").Value = "  +5.57%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextCell "D31" "4.089"
$ws.Range("E31").Value = "  -1.11%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -1.89%  "

# Row 33 - Toncoin
$ws.Range("E33").Value = "  -0.25%  "

# Row 34 - LidoDAOToken
Set-TextCell "D34" "1.846"
$ws.Range("E34").Value = "  +0.48%  "

# Row 35 - ImmutableX
Set-TextCell "D35" "0.7265"
$ws.Range("E35").Value = "  -2.13%  "

# Row 36 - ARBITRUM
Set-TextCell "D36" "1.139"
$ws.Range("E36").Value = "  -0.14%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  -2.00%  "

# Row 39 - Maker
Set-TextCell "D39" "1.220.13"
$ws.Range("E39").Value = "  +0.33%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -1.57%  "

# Row 41 - FraxShare
Set-TextCell "D41" "6.239"
$ws.Range("E41").Value = "  -3.64%  "

# Row 42 - TrustWalletToken
$ws.Range("E42").Value = "  +1.06%  "

# Row 43 - PaxDollar
Set-TextCell "D43" "1.000"
$ws.Range("E43").Value = "  -0.10%  "

# Row 44 - RocketPoolETH
Set-TextCell "D44" "1.990.68"
$ws.Range("E44").Value = "  +0.22%  "

# Row 45 - Quant
Set-TextCell "D45" "101.70"
$ws.Range("E45").Value = "  -0.21%  "

# Row 46 - Aave
Set-TextCell "D46" "65.37"
$ws.Range("E46").Value = "  -0.30%  "

# Row 47 - Mantle
Set-TextCell "D47" "0.5038"
$ws.Range("E47").Value = "  -1.04%  "

# Row 48 - BabyDogeCoin
$ws.Range("E48").Value = "  -4.10%  "

# Rows 49/50 - coins swapped: EnergySwap moves to row 49, TheSandbox to row 50
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D49" "9.140"
$ws.Range("E49").Value = "  +0.42%  "

$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell "D50" "0.4024"
$ws.Range("E50").Value = "  -1.09%  "

# Row 51 - Algorand
Set-TextCell "D51" "0.1130"
$ws.Range("E51").Value = "  +2.55%  "
